$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Determine last used row in the sheet (falls back to UsedRange if needed).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

# Column C holds the "Förändrad" (last changed) date, stored as a date serial.
# The automatic update bumped every populated row's value from 45202 to 45203
# (i.e. the "last changed" date moved forward by one day), leaving formatting
# and every other column untouched.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
